# Update countries & provincias Spain
#
# Refresh of the COVID-19 "Ciudades" data dump:
#  - three provinces near the top of the table were re-ordered
#    (Toledo / Albacete / La Rioja rotate up one slot each, picking up
#    the next reporting period's figures along the way)
#  - Guadalajara and Tenerife swap places (with updated figures)
#  - La Gomera and Arroyo de la Luz swap places (with updated figures)
#  - several rows got refreshed totals/active/recovered/deaths numbers
#  - the "last updated" timestamp moved from 16:52 to 17:22

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer
$ws.Range("A1").Value = "Datos actualizados a 18 de Abril de 2020 a las 17:22"

# Rows 11-13: Toledo / Albacete / La Rioja rotate up, new figures
$ws.Range("A11").Value = "Toledo"
$ws.Range("B11").Value = 3831
$ws.Range("C11").Value = 3838
$ws.Range("D11").Value = 10545
$ws.Range("E11").Value = 484

$ws.Range("A12").Value = "Albacete"
$ws.Range("B12").Value = 3709
$ws.Range("E12").Value = 354

$ws.Range("A13").Value = "La Rioja"
$ws.Range("B13").Value = 3648
$ws.Range("C13").Value = 2201
$ws.Range("D13").Value = 1175
$ws.Range("E13").Value = 272

# Rows 33-34: Guadalajara / Tenerife swap, new figures
$ws.Range("A33").Value = "Guadalajara"
$ws.Range("B33").Value = 1345
$ws.Range("C33").Value = 3838
$ws.Range("D33").Value = 10545
$ws.Range("E33").Value = 180

$ws.Range("A34").Value = "Tenerife"
$ws.Range("B34").Value = 1296
$ws.Range("C34").Value = 449
$ws.Range("D34").Value = 802
$ws.Range("E34").Value = 78

# Row 50: Gran Canaria refreshed figures
$ws.Range("B50").Value = 485
$ws.Range("C50").Value = 237
$ws.Range("E50").Value = 31

# Row 56: La Palma refreshed figures
$ws.Range("B56").Value = 84
$ws.Range("C56").Value = 25
$ws.Range("E56").Value = 5

# Row 57: Lanzarote refreshed figures
$ws.Range("B57").Value = 78
$ws.Range("C57").Value = 22
$ws.Range("E57").Value = 3

# Row 59: Fuerteventura refreshed figures
$ws.Range("B59").Value = 44

# Rows 62-63: La Gomera / Arroyo de la Luz swap, new figures
$ws.Range("A62").Value = "La Gomera"
$ws.Range("B62").Value = 10
$ws.Range("C62").Value = 7
$ws.Range("D62").Value = 2

$ws.Range("A63").Value = "Arroyo de la Luz"
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 7

# Row 64: El Hierro refreshed figures
$ws.Range("B64").Value = 4
$ws.Range("C64").Value = 2
